$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A (menu labels) gets a new entry "Tipos de documento" inserted at
# row 9, pushing rows 9-24 down to 10-25. Column C (descriptions) is
# untouched by this insert. Write the new column-A values directly,
# bottom-to-top so we never clobber a value before it has been moved.
$colA = @{
  25 = "Devoluciòn Proveedores nueva"
  24 = "ReKardex"
  23 = "Kardex"
  22 = "Ventas con filtro"
  21 = "Venta nueva"
  20 = "Compras con filtros"
  19 = "Compra nueva"
  18 = "Cambio de Clave"
  17 = "Cambio de Usuario"
  16 = "Búsqueda de Productos"
  15 = "Productos"
  14 = "Búsqueda de Proveedores"
  13 = "Proveedores"
  12 = "Búsqueda de Clientes"
  11 = "Validaciòn de Documento Unico"
  10 = "Clientes"
  9  = "Tipos de documento"
}

# Row 25 is brand new - give it the same plain (non-header) formatting as
# the rest of the list before writing its value.
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

foreach ($r in 25..9) {
  $ws.Cells.Item($r, 1).Value = $colA[$r]
}

# Rows 5-9 in column A now form a highlighted header-style block (same
# look as the existing A1:A4 group), matching the s="2" style change.
$ws.Range("A1").Copy()
$ws.Range("A5:A9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("C14").Select()
